# BibliotecaBaseDatos.xlsx - "Agregada interaccion basica a traves de controladores"
#
# Summary of the change:
#  - Usuarios: drop the "LibrosPrestados" column (D) entirely; rename the
#    user in row 12 (IdUsuario=10) from "Usuario_10" to "Miguel Centurion";
#    make this sheet the active tab with selection on M7.
#  - HistorialBiblioteca: append new loan/return history rows (Juan
#    borrowing "Libro 113", and Usuario_10 / Miguel Centurion borrowing &
#    returning the new book "Principito"); this sheet stops being the
#    active tab.
#  - Libros: "Libro 113" becomes "No Disponible" (it was just lent out) and
#    a brand new book row is added: 117 / Principito / Antoine de
#    Saint-Exupéry / No Disponible.

$wb = $excel.ActiveWorkbook

$wsUsuarios = $wb.Worksheets.Item("Usuarios")
$wsHistorial = $wb.Worksheets.Item("HistorialBiblioteca")
$wsLibros = $wb.Worksheets.Item("Libros")

# --- Usuarios --------------------------------------------------------------

# Remove the "LibrosPrestados" column data (keep the header cell's style,
# but clear its shared-string value; fully clear the rest of column D).
$wsUsuarios.Range("D1").ClearContents()
$wsUsuarios.Range("D2:D14").Clear()

# Row 12 (IdUsuario = 10): rename the user.
$wsUsuarios.Range("B12").Value = "Miguel Centurion"

# --- HistorialBiblioteca ----------------------------------------------------

# Juan borrows "Libro 113".
$wsHistorial.Range("A8").Value = 1
$wsHistorial.Range("B8").Value = "Juan"
$wsHistorial.Range("C8").Value = "Prestamo"
$wsHistorial.Range("D8").Value = 113
$wsHistorial.Range("E8").Value = "Libro 113"

# Usuario_10 borrows and returns the new book "Principito".
$wsHistorial.Range("A9").Value = 10
$wsHistorial.Range("B9").Value = "Usuario_10"
$wsHistorial.Range("C9").Value = "Prestamo"
$wsHistorial.Range("D9").Value = 117
$wsHistorial.Range("E9").Value = "Principito"

$wsHistorial.Range("A10").Value = 10
$wsHistorial.Range("B10").Value = "Usuario_10"
$wsHistorial.Range("C10").Value = "Devolucion"
$wsHistorial.Range("D10").Value = 117
$wsHistorial.Range("E10").Value = "Principito"

# Same user, now renamed to Miguel Centurion, borrows/returns it again.
$wsHistorial.Range("A11").Value = 10
$wsHistorial.Range("B11").Value = "Miguel Centurion"
$wsHistorial.Range("C11").Value = "Prestamo"
$wsHistorial.Range("D11").Value = 117
$wsHistorial.Range("E11").Value = "Principito"

$wsHistorial.Range("A12").Value = 10
$wsHistorial.Range("B12").Value = "Miguel Centurion"
$wsHistorial.Range("C12").Value = "Devolucion"
$wsHistorial.Range("D12").Value = 117
$wsHistorial.Range("E12").Value = "Principito"

$wsHistorial.Range("A13").Value = 10
$wsHistorial.Range("B13").Value = "Miguel Centurion"
$wsHistorial.Range("C13").Value = "Prestamo"
$wsHistorial.Range("D13").Value = 117
$wsHistorial.Range("E13").Value = "Principito"

# --- Libros -----------------------------------------------------------------

# "Libro 113" was just lent out -> no longer available.
$wsLibros.Range("D15").Value = "No Disponible"

# New book: Principito.
$wsLibros.Range("A18").Value = 117
$wsLibros.Range("B18").Value = "Principito"
$wsLibros.Range("C18").Value = "Antoine de Saint-Exupéry"
$wsLibros.Range("D18").Value = "No Disponible"

# --- Active tab / selections -------------------------------------------------
# Usuarios becomes the active sheet (activeTab=0) with selection on M7;
# HistorialBiblioteca (previously active) loses tabSelected automatically.
$wsUsuarios.Activate()
$wsUsuarios.Range("M7").Select()
